# poloniex history.xlsx -- "updated in 15 minutes"
# Row 60 (an existing "IN PROGRESS" XRP sell) is finalized as DONE, with a
# finalized date, fee, profit and duration filled in. A brand new row 61 is
# appended describing the follow-up "IN PROGRESS" XRP buy that was placed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 60: mark the trade as finished and fill in the closing details
# ---------------------------------------------------------------------
$ws.Range("H60").Value = "DONE"

$i60 = $ws.Range("I60")
$i60.Value = 42873.451099537036

$ws.Range("J60").Value = "0.06062057 USDT (0.15%)"

$k60 = $ws.Range("K60")
$k60.Value = "     ~8.5%"
$k60.Characters(6, 5).Font.Color = 5287936

$ws.Range("L60").Value = "1 day"

# ---------------------------------------------------------------------
# Row 61: brand new trade row (same timestamp as the I60 finalized date)
# ---------------------------------------------------------------------
$a61 = $ws.Range("A61")
$a61.Value = 42873.451099537036
$a61.NumberFormat = "m/d/yy h:mm"
$a61.WrapText = $true

$ws.Range("C61").Value = "        XRP"

# NOTE: the order in which brand-new shared strings are first written
# determines the index they are assigned, so G61 is written before D61,
# which is written before F61, to line up with the target workbook.
$ws.Range("G61").Value = " XRP/USDT0000008"

$d61 = $ws.Range("D61")
$d61.NumberFormat = "@"
$d61.Value = "              0.37520001`r`n`r`n"
$ws.Range("D60").Copy()
$d61.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F61").Value = "         121 XRP"
$ws.Range("E61").Value = "         0.335  USDT"
$ws.Range("H61").Value = "IN PROGRESS"

$i61 = $ws.Range("I61")
$i61.NumberFormat = "m/d/yy h:mm"
$i61.WrapText = $true

$ws.Range("K61").Value = "     "

# B61 last, so the other newly-created shared strings keep the indices
# that line up with the rest of the workbook.
$b61 = $ws.Range("B61")
$b61.Value = "            Buy"
$b61.Characters(13, 3).Font.Color = 5287936

$ws.Rows.Item(61).RowHeight = 14.25

$ws.Range("H61").Select()
